$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrigindo erro de português no campo: "Localicação CDD" -> "Localização CDD"
$ws.Range("Q1").Value = "Localização CDD"

# Atualiza a célula selecionada para refletir a edição feita pelo usuário
$ws.Range("Q2").Select()
